$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "For both hyperparameter tuning and final training, the model was trained over 40 epochs.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "For hyperparameter tuning, all models were run for 40 epochs. Final model training was run for a minimum of 5 epochs, stopping either after 40 epochs or when the validation error for one epoch exceeded that of the previous.",
    2
)
